$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 128.5
$ws.Range("I6").Value = 128.5
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 385.5
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -273.5
$ws.Range("N6").ClearContents()

$ws.Range("H18").Value = 441.66666
$ws.Range("I18").Value = 283.33334
$ws.Range("K18").Value = 283.33334
$ws.Range("M18").Value = 0.6666599999999789

$ws.Range("H32").Value = 448.3
$ws.Range("J32").Value = 440.2857
$ws.Range("L32").Value = 440.2857
$ws.Range("N32").Value = -1092.2857

$ws.Range("H132").Value = 56885.348
$ws.Range("I132").Value = 60671.867
$ws.Range("J132").Value = 6146
$ws.Range("K132").Value = 182015.601
$ws.Range("L132").Value = 18438
$ws.Range("M132").Value = -179485.601
$ws.Range("N132").Value = -23498

$ws.Range("H137").Value = 2326.5625
$ws.Range("I137").Value = 1487.8975
$ws.Range("J137").Value = 5960.778
$ws.Range("K137").Value = 4463.6925
$ws.Range("L137").Value = 17882.334
$ws.Range("M137").Value = -1913.6925
$ws.Range("N137").Value = -22982.334

$ws.Range("H138").Value = 2196.03
$ws.Range("I138").Value = 1103.5172
$ws.Range("J138").Value = 2642.2676
$ws.Range("K138").Value = 3310.5516
$ws.Range("L138").Value = 7926.8028
$ws.Range("M138").Value = 1829.4484
$ws.Range("N138").Value = -18206.8028

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 12999
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 12999
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 12999
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -13287

$ws.Range("H23").Value = 65054.668
$ws.Range("J23").Value = 57579
$ws.Range("L23").Value = 57579
$ws.Range("N23").Value = -58097

$ws.Range("H32").Value = 5865.25
$ws.Range("I32").Value = 2977.6028
$ws.Range("J32").Value = 13672.593
$ws.Range("K32").Value = 2977.6028
$ws.Range("L32").Value = 13672.593
$ws.Range("M32").Value = -2690.6028
$ws.Range("N32").Value = -14246.593

$ws.Range("H74").Value = 1499.7709
$ws.Range("I74").Value = 1068.561
$ws.Range("J74").Value = 4025.4285
$ws.Range("K74").Value = 1068.561
$ws.Range("L74").Value = 4025.4285
$ws.Range("M74").Value = -194.5609999999999
$ws.Range("N74").Value = -5773.4285

$ws.Range("H77").Value = 1499.7709
$ws.Range("I77").Value = 1068.561
$ws.Range("J77").Value = 4025.4285
$ws.Range("K77").Value = 5342.804999999999
$ws.Range("L77").Value = 20127.1425
$ws.Range("M77").Value = -974.8049999999994
$ws.Range("N77").Value = -28863.1425

$ws.Range("H102").Value = 1413.3334
$ws.Range("I102").Value = 1413.3334
$ws.Range("K102").Value = 1413.3334
$ws.Range("M102").Value = 208.6666

$ws.Range("H109").Value = 26033.572
$ws.Range("J109").Value = 26033.572
$ws.Range("L109").Value = 26033.572
$ws.Range("N109").Value = -28807.572

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 50000
$ws.Range("J35").Value = 50000
$ws.Range("L35").Value = 50000
$ws.Range("N35").Value = -50620

$ws.Range("H117").Value = 40000
$ws.Range("J117").Value = 40000
$ws.Range("L117").Value = 40000
$ws.Range("N117").Value = -49178

$ws.Range("H134").Value = 2641.7693
$ws.Range("I134").Value = 1375.7894
$ws.Range("J134").Value = 6078
$ws.Range("K134").Value = 4127.3682
$ws.Range("L134").Value = 18234
$ws.Range("M134").Value = -1592.3682
$ws.Range("N134").Value = -23304

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1464.5294
$ws.Range("J16").Value = 2118.1667
$ws.Range("L16").Value = 2118.1667
$ws.Range("N16").Value = -2692.1667

$ws.Range("H17").Value = 25000
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 25000
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 25000
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -25348

$ws.Range("H38").Value = 49999.5
$ws.Range("J38").Value = 49999.5
$ws.Range("L38").Value = 49999.5
$ws.Range("N38").Value = -50753.5

$ws.Range("H41").Value = 2800
$ws.Range("I41").Value = 2800
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 2800
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -2372
$ws.Range("N41").ClearContents()

$ws.Range("H46").Value = 49999.5
$ws.Range("J46").Value = 49999.5
$ws.Range("L46").Value = 49999.5
$ws.Range("N46").Value = -50421.5

$ws.Range("H50").Value = 27993.572
$ws.Range("J50").Value = 27993.572
$ws.Range("L50").Value = 27993.572
$ws.Range("N50").Value = -29243.572

$ws.Range("H51").Value = 30808
$ws.Range("J51").Value = 30808
$ws.Range("L51").Value = 30808
$ws.Range("N51").Value = -32280

$ws.Range("H59").Value = 35057.5
$ws.Range("J59").Value = 35057.5
$ws.Range("L59").Value = 35057.5
$ws.Range("N59").Value = -37347.5

$ws.Range("H60").Value = 33540.74
$ws.Range("J60").Value = 35834.176
$ws.Range("L60").Value = 35834.176
$ws.Range("N60").Value = -36856.176

$ws.Range("H61").Value = 30808
$ws.Range("J61").Value = 30808
$ws.Range("L61").Value = 30808
$ws.Range("N61").Value = -31504

$ws.Range("H99").Value = 3899.1177
$ws.Range("I99").Value = 2476.889
$ws.Range("J99").Value = 5499.125
$ws.Range("K99").Value = 2476.889
$ws.Range("L99").Value = 5499.125
$ws.Range("M99").Value = -978.8890000000001
$ws.Range("N99").Value = -8495.125

$ws.Range("H113").Value = 1464.5294
$ws.Range("J113").Value = 2118.1667
$ws.Range("L113").Value = 2118.1667
$ws.Range("N113").Value = -6458.1667

$ws.Range("H122").Value = 2120.2068
$ws.Range("I122").Value = 1557.2609
$ws.Range("K122").Value = 4671.7827
$ws.Range("M122").Value = -2221.7827

$ws.Range("H126").Value = 3899.1177
$ws.Range("I126").Value = 2476.889
$ws.Range("J126").Value = 5499.125
$ws.Range("K126").Value = 7430.667
$ws.Range("L126").Value = 16497.375
$ws.Range("M126").Value = -4960.667
$ws.Range("N126").Value = -21437.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 872.6070999999999
$ws.Range("J107").Value = 1260.5
$ws.Range("L107").Value = 3781.5
$ws.Range("N107").Value = -7621.5

$ws.Range("H113").Value = 568.5599999999999
$ws.Range("I113").Value = 627.9091
$ws.Range("K113").Value = 1883.7273
$ws.Range("M113").Value = 286.2727

$ws.Range("H138").Value = 3666.6667
$ws.Range("I138").Value = 3666.6667
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 11000.0001
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -5860.000100000001
$ws.Range("N138").ClearContents()

$ws.Range("H140").Value = 30633.334
$ws.Range("I140").Value = 127350
$ws.Range("K140").Value = 382050
$ws.Range("M140").Value = -376870

$ws.Range("H141").Value = 7951.1875
$ws.Range("I141").Value = 7291
$ws.Range("K141").Value = 21873
$ws.Range("M141").Value = -16693

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()

$ws.Range("H126").Value = 3945.625
$ws.Range("I126").Value = 2849.0566
$ws.Range("J126").Value = 5297.2095
$ws.Range("K126").Value = 8547.1698
$ws.Range("L126").Value = 15891.6285
$ws.Range("M126").Value = -6077.1698
$ws.Range("N126").Value = -20831.6285

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5650.2144
$ws.Range("I7").Value = 4675
$ws.Range("J7").Value = 6040.3
$ws.Range("K7").Value = 4675
$ws.Range("L7").Value = 6040.3
$ws.Range("M7").Value = -4563
$ws.Range("N7").Value = -6264.3

$ws.Range("H22").Value = 2575.05
$ws.Range("I22").Value = 2220
$ws.Range("J22").Value = 2930.1
$ws.Range("K22").Value = 2220
$ws.Range("L22").Value = 2930.1
$ws.Range("M22").Value = -1925
$ws.Range("N22").Value = -3520.1

$ws.Range("H27").Value = 2575.05
$ws.Range("I27").Value = 2220
$ws.Range("J27").Value = 2930.1
$ws.Range("K27").Value = 2220
$ws.Range("L27").Value = 2930.1
$ws.Range("M27").Value = -2113
$ws.Range("N27").Value = -3144.1

$ws.Range("H68").Value = 897.93024
$ws.Range("J68").Value = 2495
$ws.Range("L68").Value = 2495
$ws.Range("N68").Value = -3993

$ws.Range("H71").Value = 897.93024
$ws.Range("J71").Value = 2495
$ws.Range("L71").Value = 12475
$ws.Range("N71").Value = -19963

$ws.Range("H122").Value = 4203.3794
$ws.Range("I122").Value = 3133.25
$ws.Range("J122").Value = 9340
$ws.Range("K122").Value = 9399.75
$ws.Range("L122").Value = 28020
$ws.Range("M122").Value = -6949.75
$ws.Range("N122").Value = -32920

$ws.Range("H126").Value = 5650.2144
$ws.Range("I126").Value = 4675
$ws.Range("J126").Value = 6040.3
$ws.Range("K126").Value = 14025
$ws.Range("L126").Value = 18120.9
$ws.Range("M126").Value = -11555
$ws.Range("N126").Value = -23060.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 76801.5
$ws.Range("J46").Value = 76801.5
$ws.Range("L46").Value = 76801.5
$ws.Range("N46").Value = -77263.5

$ws.Range("H132").Value = 5558601.5
$ws.Range("I132").Value = 3049.8975
$ws.Range("K132").Value = 9149.692500000001
$ws.Range("M132").Value = -6619.692500000001

$ws.Range("H134").Value = 76801.5
$ws.Range("J134").Value = 76801.5
$ws.Range("L134").Value = 230404.5
$ws.Range("N134").Value = -235474.5
